$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (x)
$ws.Range("I3").Value = 0.725
$ws.Range("J3").Value = 0.725
$ws.Range("K3").Value = 0.184
$ws.Range("L3").Value = 0.184

# Row 4 (y)
$ws.Range("I4").Value = 0.181
$ws.Range("J4").Value = -0.019
$ws.Range("K4").Value = -0.019
$ws.Range("L4").Value = 0.181

# Row 5 (z)
$ws.Range("I5").Value = -0.004
$ws.Range("J5").Value = -0.004
$ws.Range("K5").Value = -0.004
$ws.Range("L5").Value = -0.004

# Row 6 (Roll)
$ws.Range("I6").Value = 3.12
$ws.Range("J6").Value = 3.12
$ws.Range("K6").Value = 3.12
$ws.Range("L6").Value = 3.12

# Row 7 (Pitch)
$ws.Range("I7").Value = 1.444
$ws.Range("J7").Value = 1.44
$ws.Range("K7").Value = 1.44
$ws.Range("L7").Value = 1.44

# Row 8 (Yaw)
$ws.Range("I8").Value = 3.032
$ws.Range("J8").Value = 3.032
$ws.Range("K8").Value = 3.032
$ws.Range("L8").Value = 3.032

# Update the selected cell/range in the sheet view to L5
$ws.Range("L5").Select()
